# Chap Research 2015 Kickoff Presentation - Jan 2015 update
#
# Semantic edit (confirmed by reverse-engineering the target OOXML diff):
#   1. Delete slide "What It Takes" (slide 9)
#   2. Delete slide "Space to Meet" (slide 10, which was slide 10 before any
#      deletion -- delete the higher index first so indices stay stable)
#   3. On the "The Process" slide (slide 8, title shape), change the title
#      text to "Our Process" (typed as two runs: "Our " then "Process",
#      matching how PowerPoint records a retyped/reworded run boundary)
#
# All of the other differences in the canonical OOXML (renumbered
# relationship ids, notesMasterId, cached slide-number fields, cached
# datetimeFigureOut field text, etc.) are automatic side effects that
# PowerPoint itself recomputes on save after the slide count/ordering
# changes above -- they are not separate user actions.

$p = $ppt.ActivePresentation

# --- 1 & 2: remove the two slides -----------------------------------------
# Delete slide 10 ("Space to Meet") before slide 9 ("What It Takes") so the
# index of slide 9 doesn't shift out from under us.
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()

# --- 3: retitle "The Process" slide to "Our Process" ----------------------
$titleRange = $p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Our "
$titleRange.InsertAfter("Process")
